$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 788.1599833333333
$ws.Range("N2").Value = 2364.47995
$ws.Range("O2").Value = 0.8397951873720987
$ws.Range("P2").Value = 0.8397951873720988
$ws.Range("Q2").Value = 133920.7213891924
$ws.Range("R2").Value = 1205286.492502732
$ws.Range("S2").Value = 0.3729685731900123
$ws.Range("T2").Value = 0.3729685731900124
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.01890163353898316
$ws.Range("P3").Value = 0.01890163353898317
$ws.Range("Q3").Value = 3014.211604255354
$ws.Range("R3").Value = 27127.90443829818
$ws.Range("S3").Value = 0.008394565005849964
$ws.Range("T3").Value = 0.008394565005849968
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 131.4690986666667
$ws.Range("N4").Value = 394.407296
$ws.Range("O4").Value = 0.1400821136357036
$ws.Range("P4").Value = 0.1400821136357036
$ws.Range("Q4").Value = 22338.65827514449
$ws.Range("R4").Value = 201047.9244763004
$ws.Range("S4").Value = 0.0622130572284408
$ws.Range("T4").Value = 0.06221305722844081
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("M5").Value = 1.145987666666667
$ws.Range("N5").Value = 3.437963
$ws.Range("O5").Value = 0.001221065453214498
$ws.Range("P5").Value = 0.001221065453214498
$ws.Range("Q5").Value = 194.7212472955637
$ws.Range("R5").Value = 1752.491225660073
$ws.Range("S5").Value = 0.0005422977491477795
$ws.Range("T5").Value = 0.0005422977491477797
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 788.1599833333333
$ws.Range("N6").Value = 2364.47995
$ws.Range("O6").Value = 0.8397951873720987
$ws.Range("P6").Value = 0.8397951873720988
$ws.Range("Q6").Value = 53896.14750317595
$ws.Range("R6").Value = 485065.3275285835
$ws.Range("S6").Value = 0.1501005148880583
$ws.Range("T6").Value = 0.1501005148880583
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.01890163353898316
$ws.Range("P7").Value = 0.01890163353898317
$ws.Range("S7").Value = 0.00337837721516934
$ws.Range("T7").Value = 0.00337837721516934
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("M8").Value = 131.4690986666667
$ws.Range("N8").Value = 394.407296
$ws.Range("O8").Value = 0.1400821136357036
$ws.Range("P8").Value = 0.1400821136357036
$ws.Range("Q8").Value = 8990.151852014977
$ws.Range("R8").Value = 80911.36666813478
$ws.Range("S8").Value = 0.02503753022105635
$ws.Range("T8").Value = 0.02503753022105635
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("M9").Value = 1.145987666666667
$ws.Range("N9").Value = 3.437963
$ws.Range("O9").Value = 0.001221065453214498
$ws.Range("P9").Value = 0.001221065453214498
$ws.Range("Q9").Value = 78.36520709700299
$ws.Range("R9").Value = 705.286863873027
$ws.Range("S9").Value = 0.0002182467296735139
$ws.Range("T9").Value = 0.0002182467296735139
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 788.1599833333333
$ws.Range("N10").Value = 2364.47995
$ws.Range("O10").Value = 0.8397951873720987
$ws.Range("P10").Value = 0.8397951873720988
$ws.Range("Q10").Value = 41988.93464552943
$ws.Range("R10").Value = 377900.4118097648
$ws.Range("S10").Value = 0.1169389836170313
$ws.Range("T10").Value = 0.1169389836170314
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.01890163353898316
$ws.Range("P11").Value = 0.01890163353898317
$ws.Range("Q11").Value = 945.0631145501605
$ws.Range("R11").Value = 8505.568030951445
$ws.Range("S11").Value = 0.00263199628669808
$ws.Range("T11").Value = 0.002631996286698081
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("M12").Value = 131.4690986666667
$ws.Range("N12").Value = 394.407296
$ws.Range("O12").Value = 0.1400821136357036
$ws.Range("P12").Value = 0.1400821136357036
$ws.Range("Q12").Value = 7003.968113776554
$ws.Range("R12").Value = 63035.71302398898
$ws.Range("S12").Value = 0.01950601794080834
$ws.Range("T12").Value = 0.01950601794080835
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("M13").Value = 1.145987666666667
$ws.Range("N13").Value = 3.437963
$ws.Range("O13").Value = 0.001221065453214498
$ws.Range("P13").Value = 0.001221065453214498
$ws.Range("Q13").Value = 61.05207351018066
$ws.Range("R13").Value = 549.4686615916258
$ws.Range("S13").Value = 0.0001700297348399845
$ws.Range("T13").Value = 0.0001700297348399846
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 788.1599833333333
$ws.Range("N14").Value = 2364.47995
$ws.Range("O14").Value = 0.8397951873720987
$ws.Range("P14").Value = 0.8397951873720988
$ws.Range("Q14").Value = 71736.96814958862
$ws.Range("R14").Value = 645632.7133462975
$ws.Range("S14").Value = 0.1997871156769967
$ws.Range("T14").Value = 0.1997871156769968
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.01890163353898316
$ws.Range("P15").Value = 0.01890163353898317
$ws.Range("Q15").Value = 1614.614972258034
$ws.Range("R15").Value = 14531.5347503223
$ws.Range("S15").Value = 0.004496695031265779
$ws.Range("T15").Value = 0.00449669503126578
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("M16").Value = 131.4690986666667
$ws.Range("N16").Value = 394.407296
$ws.Range("O16").Value = 0.1400821136357036
$ws.Range("P16").Value = 0.1400821136357036
$ws.Range("Q16").Value = 11966.09158437456
$ws.Range("R16").Value = 107694.824259371
$ws.Range("S16").Value = 0.03332550824539809
$ws.Range("T16").Value = 0.0333255082453981
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("M17").Value = 1.145987666666667
$ws.Range("N17").Value = 3.437963
$ws.Range("O17").Value = 0.001221065453214498
$ws.Range("P17").Value = 0.001221065453214498
$ws.Range("Q17").Value = 104.3058293771804
$ws.Range("R17").Value = 938.752464394624
$ws.Range("S17").Value = 0.0002904912395532195
$ws.Range("T17").Value = 0.0002904912395532196

Write-Host "Applied 174 cell updates"
